# Apply CodeBlockData updates: new BlockName rows (Up/Down/Left/Right/Attack*/Condition)
# plus column width / selection adjustments for the new BlockName column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Up"),
    @(2, "Down"),
    @(3, "Left"),
    @(4, "Right"),
    @(5, "AttackFire"),
    @(6, "AttackWater"),
    @(7, "AttackGrass"),
    @(8, "Condition")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

$ws.Columns.Item(1).ColumnWidth = 15.375
$ws.Columns.Item(2).ColumnWidth = 23.25

$ws.Range("A6").Select()
